# Rename three sheets to insert a hyphen before the trailing "RI"
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("TGZ-S-48-50_100RI").Name = "TGZ-S-48-50_100-RI"
$wb.Worksheets.Item("TGZ-S-48-100_250RI").Name = "TGZ-S-48-100_250-RI"
$wb.Worksheets.Item("TGZ-S-48-100_300RI").Name = "TGZ-S-48-100_300-RI"

# Update the selected/active cell remembered in each sheet's view
$ws1 = $wb.Worksheets.Item("TGZ-S-48-50_100-RI")
$ws1.Activate()
$ws1.Range("C39").Select()

$ws2 = $wb.Worksheets.Item("TGZ-S-48-100_300-RI")
$ws2.Activate()
$ws2.Range("F20").Select()
